# Apply the data-row updates described by the diff for rows 25-28 on the
# active worksheet. The content of rows 25-28 is effectively rotated
# (row25<-old row28, row26<-old row27, row27<-old row26, row28<-old row25)
# with a new "Taxonsorteringsordning" (column B) value on each row, plus the
# corresponding additions/removals of the optional comment/biotope cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 25 ----
$ws.Range("A25").Value = 112306119
$ws.Range("B25").Value = 90826
$ws.Range("D25").Value = "LC"
$ws.Range("E25").Value = 4366
$ws.Range("F25").Value = "Skarp dropptaggsvamp"
$ws.Range("G25").Value = "Hydnellum peckii"
$ws.Range("H25").Value = "Banker"
$ws.Range("AC25").ClearContents()
$ws.Range("AH25").ClearContents()
$ws.Range("AI25").ClearContents()

# ---- Row 26 ----
$ws.Range("A26").Value = 112306159
$ws.Range("B26").Value = 90810
$ws.Range("D26").Value = "LC"
$ws.Range("E26").Value = 4363
$ws.Range("F26").Value = "Zontaggsvamp"
$ws.Range("G26").Value = "Hydnellum concrescens"
$ws.Range("H26").Value = "(Pers.) Banker"
$ws.Range("J26").ClearContents()
$ws.Range("AC26").Value = "Efter stigen i början nära vändplatsen"
$ws.Range("AH26").ClearContents()
$ws.Range("AI26").ClearContents()

# ---- Row 27 ----
$ws.Range("A27").Value = 112306179
$ws.Range("B27").Value = 89114
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 5754
$ws.Range("F27").Value = "Gultoppig fingersvamp"
$ws.Range("G27").Value = "Ramaria testaceoflava"
$ws.Range("H27").Value = "(Bres.) Corner"
$ws.Range("J27").Value = "fruktkroppar"
$ws.Range("AC27").Value = "Där stigen delar sig ned mot myren"
$ws.Range("AH27").Value = "Skogsmark"
$ws.Range("AI27").Value = "Barrblandskog kalkpåverkad."

# ---- Row 28 ----
$ws.Range("A28").Value = 112306136
$ws.Range("B28").Value = 90832
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 4368
$ws.Range("F28").Value = "Dofttaggsvamp"
$ws.Range("G28").Value = "Hydnellum suaveolens"
$ws.Range("H28").Value = "(Scop.:Fr.) P. Karst."
$ws.Range("AC28").Value = "Förekommer på flera platser"
$ws.Range("AH28").Value = "Skogsmark"
$ws.Range("AI28").Value = "Barrblandskog kalkpåverkad."
